# Updated return statistics + new scenario generation method
# Applies new forecast values (column B, "MSTL") for rows 2-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 106.197998046875
    3  = 106.1641540527344
    4  = 106.2518692016602
    5  = 106.2052459716797
    6  = 104.5456314086914
    7  = 104.4966812133789
    8  = 103.8380813598633
    9  = 103.7358627319336
    10 = 104.6287536621094
    11 = 104.5373153686523
    12 = 109.8541107177734
    13 = 109.7744293212891
    14 = 123.1528549194336
    15 = 123.0840606689453
    16 = 140.2832336425781
    17 = 140.2240447998047
    18 = 151.8777313232422
    19 = 151.8302764892578
    20 = 162.6111602783203
    21 = 162.5661926269531
    22 = 162.53466796875
    23 = 162.4917144775391
    24 = 156.7200622558594
    25 = 156.6782531738281
    26 = 159.3792266845703
    27 = 159.3375854492188
    28 = 154.4311981201172
    29 = 154.3886871337891
    30 = 155.5289001464844
    31 = 155.4848327636719
    32 = 164.8945922851562
    33 = 164.8486328125
    34 = 187.7915954589844
    35 = 187.7430725097656
    36 = 211.5
    37 = 211.4480743408203
    38 = 186.4666290283203
    39 = 186.4101715087891
    40 = 163.5970153808594
    41 = 163.5346374511719
    42 = 147.7666931152344
    43 = 147.697021484375
    44 = 136.4794616699219
    45 = 136.4012603759766
    46 = 125.8640823364258
    47 = 125.7766876220703
    48 = 124.4602890014648
    49 = 124.3642807006836
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
